$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.821.13'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '2.083.24'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'233.80"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = "'58.72"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'0.394"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('D10').Value = "'0.0788"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('E11').Value = '  +3.19%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').Value = "'14.98"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.98%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.390.26'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('D14').Value = "'21.28"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').Value = "'0.785"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').Value = "'5.37"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.70%  '
$ws.Range('D17').Value = '2.093.19'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Value = '37.735.20'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D20').Value = "'71.55"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').Value = '0.0₃0842'
$ws.Range('E21').Value = '  +1.65%  '
$ws.Range('D22').Value = "'230.08"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('E25').Value = '  +1.37%  '
$ws.Range('D26').Value = "'10.03"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.26%  '
$ws.Range('D27').Value = "'172.13"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.15%  '
$ws.Range('E28').Value = '  -1.79%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value = "'1.41"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = "'19.53"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('D32').Value = "'4.75"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('D33').Value = "'0.0635"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('E34').Value = '  -0.83%  '
$ws.Range('E35').Value = '  -1.84%  '
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = "'5.44"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('D40').Value = "'0.0235"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.55%  '
$ws.Range('D41').Value = "'101.42"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.93%  '
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('E43').Value = '  -1.02%  '
$ws.Range('D44').Value = "'16.95"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.86%  '
$ws.Range('D45').Value = '1.448.48'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = "'4.10"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.76%  '
$ws.Range('D49').Value = "'7.35"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').Value = "'3.00"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('D51').Value = '2.275.36'
$ws.Range('E51').Value = '  -0.01%  '
